# Atualizacao de bases das ligas (Germany Landesliga), 02-04-2024.
# A newly-scraped batch of match results caused several existing match
# rows to be reordered/updated; apply the resulting per-row field changes
# in place (columns B and F..AC), leaving A (row index), C, D and E
# (Div / Div Original Name / Date) untouched for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (id=2) now holds the match previously stored in row 5
$ws.Cells.Item(4, 2).Value() = 6781315
$ws.Cells.Item(4, 6).Value() = "SSV Markranstadt"
$ws.Cells.Item(4, 7).Value() = "BSC Rapid Chemnitz"
$ws.Cells.Item(4, 8).Value() = 2
$ws.Cells.Item(4, 9).Value() = 0
$ws.Cells.Item(4, 10).Value() = "H"
$ws.Cells.Item(4, 11).Value() = 1.25
$ws.Cells.Item(4, 12).Value() = 4.75
$ws.Cells.Item(4, 13).Value() = 10
$ws.Cells.Item(4, 14).Value() = 1.222
$ws.Cells.Item(4, 15).Value() = 5.25
$ws.Cells.Item(4, 16).Value() = 8.5
$ws.Cells.Item(4, 17).Value() = -2
$ws.Cells.Item(4, 18).Value() = 1.925
$ws.Cells.Item(4, 19).Value() = 1.875
$ws.Cells.Item(4, 20).Value() = 3.5
$ws.Cells.Item(4, 21).Value() = 1.775
$ws.Cells.Item(4, 22).Value() = 1.925
$ws.Cells.Item(4, 23).Value() = 0.222
$ws.Cells.Item(4, 24).Value() = -1
$ws.Cells.Item(4, 25).Value() = -1
$ws.Cells.Item(4, 26).Value() = 0
$ws.Cells.Item(4, 27).Value() = -0
$ws.Cells.Item(4, 28).Value() = -1
$ws.Cells.Item(4, 29).Value() = 0.925

# Row 5 (id=3) now holds the match previously stored in row 4
$ws.Cells.Item(5, 2).Value() = 6781316
$ws.Cells.Item(5, 6).Value() = "SV Schott Jena"
$ws.Cells.Item(5, 7).Value() = "SV 09 Arnstadt"
$ws.Cells.Item(5, 8).Value() = 0
$ws.Cells.Item(5, 9).Value() = 2
$ws.Cells.Item(5, 10).Value() = "A"
$ws.Cells.Item(5, 11).Value() = 5
$ws.Cells.Item(5, 12).Value() = 4.5
$ws.Cells.Item(5, 13).Value() = 1.45
$ws.Cells.Item(5, 14).Value() = 6.5
$ws.Cells.Item(5, 15).Value() = 4.333
$ws.Cells.Item(5, 16).Value() = 1.363
$ws.Cells.Item(5, 17).Value() = 1.5
$ws.Cells.Item(5, 18).Value() = 1.825
$ws.Cells.Item(5, 19).Value() = 1.975
$ws.Cells.Item(5, 20).Value() = 3
$ws.Cells.Item(5, 21).Value() = 1.825
$ws.Cells.Item(5, 22).Value() = 1.975
$ws.Cells.Item(5, 23).Value() = -1
$ws.Cells.Item(5, 24).Value() = -1
$ws.Cells.Item(5, 25).Value() = 0.363
$ws.Cells.Item(5, 26).Value() = -1
$ws.Cells.Item(5, 27).Value() = 0.9750000000000001
$ws.Cells.Item(5, 28).Value() = -1
$ws.Cells.Item(5, 29).Value() = 0.9750000000000001

# Row 11 (id=9) now holds the match previously stored in row 12
$ws.Cells.Item(11, 2).Value() = 7035047
$ws.Cells.Item(11, 6).Value() = "SC Dsseldorf West"
$ws.Cells.Item(11, 7).Value() = "VfL Viktoria JuchenGarzweiler"
$ws.Cells.Item(11, 8).Value() = 3
$ws.Cells.Item(11, 9).Value() = 4
$ws.Cells.Item(11, 10).Value() = "A"
$ws.Cells.Item(11, 11).Value() = 1.909
$ws.Cells.Item(11, 12).Value() = 3.75
$ws.Cells.Item(11, 13).Value() = 3.1
$ws.Cells.Item(11, 14).Value() = 2.2
$ws.Cells.Item(11, 15).Value() = 3.6
$ws.Cells.Item(11, 16).Value() = 2.625
$ws.Cells.Item(11, 17).Value() = -0.25
$ws.Cells.Item(11, 18).Value() = 2
$ws.Cells.Item(11, 19).Value() = 1.8
$ws.Cells.Item(11, 20).Value() = 3
$ws.Cells.Item(11, 21).Value() = 1.825
$ws.Cells.Item(11, 22).Value() = 1.975
$ws.Cells.Item(11, 23).Value() = -1
$ws.Cells.Item(11, 24).Value() = -1
$ws.Cells.Item(11, 25).Value() = 1.625
$ws.Cells.Item(11, 26).Value() = -1
$ws.Cells.Item(11, 27).Value() = 0.8
$ws.Cells.Item(11, 28).Value() = 0.825
$ws.Cells.Item(11, 29).Value() = -1

# Row 12 (id=10) now holds the match previously stored in row 13
$ws.Cells.Item(12, 2).Value() = 7035046
$ws.Cells.Item(12, 6).Value() = "Cronenberger SC"
$ws.Cells.Item(12, 7).Value() = "FC Viersen"
$ws.Cells.Item(12, 8).Value() = 0
$ws.Cells.Item(12, 9).Value() = 2
$ws.Cells.Item(12, 10).Value() = "A"
$ws.Cells.Item(12, 11).Value() = 2
$ws.Cells.Item(12, 12).Value() = 3.6
$ws.Cells.Item(12, 13).Value() = 3
$ws.Cells.Item(12, 14).Value() = 2
$ws.Cells.Item(12, 15).Value() = 3.6
$ws.Cells.Item(12, 16).Value() = 3
$ws.Cells.Item(12, 17).Value() = -0.25
$ws.Cells.Item(12, 18).Value() = 1.8
$ws.Cells.Item(12, 19).Value() = 2
$ws.Cells.Item(12, 20).Value() = 2.75
$ws.Cells.Item(12, 21).Value() = 1.8
$ws.Cells.Item(12, 22).Value() = 2
$ws.Cells.Item(12, 23).Value() = -1
$ws.Cells.Item(12, 24).Value() = -1
$ws.Cells.Item(12, 25).Value() = 2
$ws.Cells.Item(12, 26).Value() = -1
$ws.Cells.Item(12, 27).Value() = 1
$ws.Cells.Item(12, 28).Value() = -1
$ws.Cells.Item(12, 29).Value() = 1

# Row 13 (id=11) now holds the match previously stored in row 11
$ws.Cells.Item(13, 2).Value() = 7035048
$ws.Cells.Item(13, 6).Value() = "SG Unterrath"
$ws.Cells.Item(13, 7).Value() = "TuRU Dsseldorf"
$ws.Cells.Item(13, 8).Value() = 1
$ws.Cells.Item(13, 9).Value() = 0
$ws.Cells.Item(13, 10).Value() = "H"
$ws.Cells.Item(13, 11).Value() = 3.25
$ws.Cells.Item(13, 12).Value() = 4
$ws.Cells.Item(13, 13).Value() = 1.8
$ws.Cells.Item(13, 14).Value() = 2.9
$ws.Cells.Item(13, 15).Value() = 4
$ws.Cells.Item(13, 16).Value() = 1.95
$ws.Cells.Item(13, 17).Value() = 0.5
$ws.Cells.Item(13, 18).Value() = 1.8
$ws.Cells.Item(13, 19).Value() = 2
$ws.Cells.Item(13, 20).Value() = 3
$ws.Cells.Item(13, 21).Value() = 1.75
$ws.Cells.Item(13, 22).Value() = 1.95
$ws.Cells.Item(13, 23).Value() = 1.9
$ws.Cells.Item(13, 24).Value() = -1
$ws.Cells.Item(13, 25).Value() = -1
$ws.Cells.Item(13, 26).Value() = 0.8
$ws.Cells.Item(13, 27).Value() = -1
$ws.Cells.Item(13, 28).Value() = -1
$ws.Cells.Item(13, 29).Value() = 0.95

# Row 46 (id=44) now holds the match previously stored in row 47
$ws.Cells.Item(46, 2).Value() = 7511941
$ws.Cells.Item(46, 6).Value() = "SpVgg SterkradeNord"
$ws.Cells.Item(46, 7).Value() = "BlauWeiss Mintard"
$ws.Cells.Item(46, 8).Value() = 3
$ws.Cells.Item(46, 9).Value() = 1
$ws.Cells.Item(46, 10).Value() = "H"
$ws.Cells.Item(46, 11).Value() = 3.4
$ws.Cells.Item(46, 12).Value() = 4.2
$ws.Cells.Item(46, 13).Value() = 1.727
$ws.Cells.Item(46, 14).Value() = 2.8
$ws.Cells.Item(46, 15).Value() = 4.2
$ws.Cells.Item(46, 16).Value() = 1.95
$ws.Cells.Item(46, 17).Value() = 0.5
$ws.Cells.Item(46, 18).Value() = 1.8
$ws.Cells.Item(46, 19).Value() = 2
$ws.Cells.Item(46, 20).Value() = 3.5
$ws.Cells.Item(46, 21).Value() = 1.875
$ws.Cells.Item(46, 22).Value() = 1.925
$ws.Cells.Item(46, 23).Value() = 1.8
$ws.Cells.Item(46, 24).Value() = -1
$ws.Cells.Item(46, 25).Value() = -1
$ws.Cells.Item(46, 26).Value() = 0.8
$ws.Cells.Item(46, 27).Value() = -1
$ws.Cells.Item(46, 28).Value() = 0.875
$ws.Cells.Item(46, 29).Value() = -1

# Row 47
$ws.Cells.Item(47, 2).Value() = 7511940
$ws.Cells.Item(47, 6).Value() = "ASV Suchteln"
$ws.Cells.Item(47, 7).Value() = "Holzheimer SG"
$ws.Cells.Item(47, 8).Value() = 3
$ws.Cells.Item(47, 9).Value() = 3
$ws.Cells.Item(47, 10).Value() = "D"
$ws.Cells.Item(47, 11).Value() = 2.75
$ws.Cells.Item(47, 12).Value() = 3.6
$ws.Cells.Item(47, 13).Value() = 2.1
$ws.Cells.Item(47, 14).Value() = 4
$ws.Cells.Item(47, 15).Value() = 4
$ws.Cells.Item(47, 16).Value() = 1.6
$ws.Cells.Item(47, 17).Value() = 1
$ws.Cells.Item(47, 18).Value() = 1.8
$ws.Cells.Item(47, 19).Value() = 2
$ws.Cells.Item(47, 20).Value() = 3.5
$ws.Cells.Item(47, 21).Value() = 1.85
$ws.Cells.Item(47, 22).Value() = 1.95
$ws.Cells.Item(47, 23).Value() = -1
$ws.Cells.Item(47, 24).Value() = 3
$ws.Cells.Item(47, 25).Value() = -1
$ws.Cells.Item(47, 26).Value() = 0.8
$ws.Cells.Item(47, 27).Value() = -1
$ws.Cells.Item(47, 28).Value() = 0.8500000000000001
$ws.Cells.Item(47, 29).Value() = -1
